# Automatic update of files.
# Cyclic shift of data rows 2-5: each row's content moves up one row
# (row3->row2, row4->row3, row5->row4), and the original row 2 content
# wraps around to row 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (was row 3's data) ---
$ws.Range("A2").Value = 94746267
$ws.Range("B2").Value = 96334
$ws.Range("D2").Value = "VU"
$ws.Range("E2").Value = 220787
$ws.Range("F2").Value = "Knärot"
$ws.Range("G2").Value = "Goodyera repens"
$ws.Range("H2").Value = "(L.) R. Br."
$ws.Range("I2").Value = "12"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "fullt utvecklade blad"
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("Q2").Value = 563670.0922861055
$ws.Range("R2").Value = 6711641.740220776
$ws.Range("AF2").Value = ""

# --- Row 3 (was row 4's data) ---
$ws.Range("A3").Value = 94746350
$ws.Range("B3").Value = 101120
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 222002
$ws.Range("F3").Value = "Underviol"
$ws.Range("G3").Value = "Viola mirabilis"
$ws.Range("H3").Value = "L."
$ws.Range("I3").Value = ""
$ws.Range("K3").Value = ""
$ws.Range("Q3").Value = 563673.5095373251
$ws.Range("R3").Value = 6711615.711973636

# --- Row 4 (was row 5's data) ---
$ws.Range("A4").Value = 94746394
$ws.Range("B4").Value = 96334
$ws.Range("D4").Value = "VU"
$ws.Range("E4").Value = 220787
$ws.Range("F4").Value = "Knärot"
$ws.Range("G4").Value = "Goodyera repens"
$ws.Range("H4").Value = "(L.) R. Br."
$ws.Range("I4").Value = "2"
$ws.Range("K4").Value = "blomknopp"
$ws.Range("Q4").Value = 563679.3894672301
$ws.Range("R4").Value = 6711617.784626649

# --- Row 5 (was row 2's data) ---
$ws.Range("A5").Value = 94747241
$ws.Range("B5").Value = 5113
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 100526
$ws.Range("F5").Value = "Bronshjon"
$ws.Range("G5").Value = "Callidium coriaceum"
$ws.Range("H5").Value = "Paykull, 1800"
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("M5").Value = "äldre gnagspår"
$ws.Range("N5").Value = ""
$ws.Range("Q5").Value = 563857.9908253724
$ws.Range("R5").Value = 6711666.716977899
$ws.Range("AF5").Value = ""
